# [NOBTS] Add the script structure ppt
#
# Cleans up leftover "retyped" runs on slide 2:
#   - The two "TestRunner<br/>Object" rectangles: the "TestRunner" run
#     carries a stray err="1" (spellcheck) flag and is immediately
#     followed by a redundant empty run. Retyping "TestRunner" removes
#     both: the empty run absorbs the retyped text (picking up its own,
#     error-free rPr) and the err-flagged run disappears.
#   - The "JythonScript Engine" rectangle: "JythonScript" (err="1") and
#     the following " " run merge into a single, error-free
#     "JythonScript " run.
#
# Both fixes are done the way PowerPoint itself arrives at this XML
# shape: clear the old run's characters, then retype the word into the
# (now empty) range immediately preceding the rest of the text, so the
# surviving run keeps the clean (err-free) rPr instead of the flagged one.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# -- "직사각형 8" (id 9) / "직사각형 9" (id 10): TestRunner<br/>Object --
# Erasing "TestRunner" (chars 1-10) leaves the text frame starting with
# the pre-existing empty run (no err) immediately followed by the <br/>;
# retyping into that now-empty leading run keeps its error-free rPr.
foreach ($idx in 6, 7) {
    $tr = $s.Shapes.Item($idx).TextFrame.TextRange
    $tr.Characters(1, 10).Text = ""
    $tr.Characters(1, 0).Text = "TestRunner"
}

# -- "직사각형 54" (id 55): JythonScript Engine --
# Erasing "JythonScript" (chars 1-12, the err-flagged run) leaves the
# text starting with the separate, error-free " " run; retyping
# "JythonScript " into that surviving single-space run merges the two
# words into one clean run while leaving "Engine" untouched.
$tr23 = $s.Shapes.Item(23).TextFrame.TextRange
$tr23.Characters(1, 12).Text = ""
$tr23.Characters(1, 1).Text = "JythonScript "
